$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the rows needed to make room for the new average rows ---
# Old row 15 (header of 2nd block) -> new row 16 (shift down by 1)
$ws.Rows.Item(14).EntireRow.Insert()
# Old row 27 (header of 3rd block, now at row 28 after first insert) -> new row 29 (shift down by 1 more, total +2)
$ws.Rows.Item(28).EntireRow.Insert()

# --- Average row after the 1st block (case1-10), new row 14 ---
$ws.Range("G14").Formula = "=AVERAGE(G4:G13)"
$ws.Range("H14").Formula = "=AVERAGE(H4:H13)"

# --- Average row after the 2nd block (case1-10), new row 27 ---
$ws.Range("G27").Formula = "=AVERAGE(G17:G26)"
$ws.Range("H27").Formula = "=AVERAGE(H17:H26)"

# --- Average row after the 3rd block (case1-10), new row 40 ---
$ws.Range("G40").Formula = "=AVERAGE(G30:G39)"
$ws.Range("H40").Formula = "=AVERAGE(H30:H39)"

# --- New 4th block: header + 3 new "민경진" rows, rows 47-50 ---
$ws.Range("D47").Value = "사용자 승리"
$ws.Range("E47").Value = "컴퓨터 승리"
$ws.Range("F47").Value = "전체 게임 수"
$ws.Range("G47").Value = "컴퓨터 승률"

$ws.Range("C48").Value = "민경진(1)"
$ws.Range("D48").Value = 29
$ws.Range("E48").Value = 39
$ws.Range("F48").Value = 100
$ws.Range("G48").Value = 39

$ws.Range("C49").Value = "민경진(2)"
$ws.Range("D49").Value = 26
$ws.Range("E49").Value = 41
$ws.Range("F49").Value = 100
$ws.Range("G49").Value = 41

$ws.Range("C50").Value = "민경진(3)"
$ws.Range("D50").Value = 25
$ws.Range("E50").Value = 52
$ws.Range("F50").Value = 100
$ws.Range("G50").Value = 52

# --- Update sheet view to match the saved selection/scroll position ---
$ws.Range("G51").Select() | Out-Null
